# This workbook is a test/status tracker for AutoHotkey functions.
# Column A = function name, column B = test status ("Pass"/"pass"/"N/a"/"n/a"/etc.),
# with an AutoFilter on column B set to show only blank ("not yet tested") rows.
#
# The edit fills in the previously-blank status for six more functions with
# "n/a" (they are not applicable / not covered by a test), which in turn
# causes those rows - plus the already-failing-filter row 4 - to be hidden
# by the existing "blanks only" AutoFilter on column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "n/a" status for the functions that did not have one yet.
$ws.Range("B40").Value = "n/a"   # Move_Message_Box()
$ws.Range("B42").Value = "n/a"   # Pausescript()
$ws.Range("B50").Value = "n/a"   # SerialbreakquestionGUI()
$ws.Range("B51").Value = "n/a"   # Serials_GUI_Screen()
$ws.Range("B55").Value = "n/a"   # UnPausescript()
$ws.Range("B56").Value = "n/a"   # Versioncheck()

# The AutoFilter on column B (blanks only) now hides these rows, since they
# are no longer blank. Row 4 already had a non-blank status ("pass") and is
# also (re-)hidden once the filter state is refreshed.
$ws.Rows.Item(4).Hidden = $true
$ws.Rows.Item(40).Hidden = $true
$ws.Rows.Item(42).Hidden = $true
$ws.Rows.Item(50).Hidden = $true
$ws.Rows.Item(51).Hidden = $true
$ws.Rows.Item(55).Hidden = $true
$ws.Rows.Item(56).Hidden = $true

# The active selection moved from B44 to A44.
$null = $ws.Range("A44").Select()
